$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rework "3. HV test" results (rows 60-61) and make room for the expanded
#    "HV test" + new "4. Visual inspection" sections.
# ---------------------------------------------------------------------------

# Update the first two instruction lines of the HV test section.
$ws.Range("B60").Value = "1) Do not connect the adapter. Turn HV on, read after 60 sec"
$ws.Range("B61").Value = "2) Connect the adapter. Turn HV on, read after 60 sec"

# Insert 7 fresh rows right after row 61 so the rest of the sheet (old rows
# 63..72) is pushed down to 70..79, and we have room for the new layout.
$ws.Rows("62:68").Insert()

# Wipe everything in the region we are about to rebuild (this also clears
# the old "3. HV test" tail / "4. Summary" block that got shifted into
# rows 70-77 by the insert above).
$ws.Range("A62:G78").Clear()

# ---------------------------------------------------------------------------
# 2) New HV-test measurement block (rows 62, 64-68)
# ---------------------------------------------------------------------------

$ws.Range("B62").Value = "Keithley should not trip any time."

$ws.Range("B64").Value = "No adapter"
$ws.Range("C64").Value = 0.043
$ws.Range("C64").NumberFormat = "0.0000"
$ws.Range("C64").HorizontalAlignment = -4108
$ws.Range("C64").VerticalAlignment = -4160
$ws.Range("D64").Value = "µA"
$ws.Range("D64").HorizontalAlignment = -4131
$ws.Range("D64").VerticalAlignment = -4160

$ws.Range("B65").Value = "With adapt."
$ws.Range("C65").Value = 0.05
$ws.Range("C65").NumberFormat = "0.0000"
$ws.Range("C65").HorizontalAlignment = -4108
$ws.Range("C65").VerticalAlignment = -4160
$ws.Range("D65").Value = "µA"
$ws.Range("D65").HorizontalAlignment = -4131
$ws.Range("D65").VerticalAlignment = -4160

$ws.Range("B66").Value = "Difference"
$ws.Range("C66").Formula = "=C65-C64"
$ws.Range("C66").NumberFormat = "0.0000"
$ws.Range("C66").HorizontalAlignment = -4108
$ws.Range("C66").VerticalAlignment = -4160
$ws.Range("D66").Value = "µA"
$ws.Range("D66").HorizontalAlignment = -4131
$ws.Range("D66").VerticalAlignment = -4160
$ws.Range("E66").Value = "Limit:"
$ws.Range("F66").Value = 0.01
$ws.Range("G66").Value = "µA"
$ws.Range("G66").HorizontalAlignment = -4131
$ws.Range("G66").VerticalAlignment = -4160

# Blank helper row under the measurement block, keeps the same formatting.
$ws.Range("C67").NumberFormat = "0.0000"
$ws.Range("C67").HorizontalAlignment = -4108
$ws.Range("C67").VerticalAlignment = -4160
$ws.Range("D67").HorizontalAlignment = -4131
$ws.Range("D67").VerticalAlignment = -4160
$ws.Range("G67").HorizontalAlignment = -4131
$ws.Range("G67").VerticalAlignment = -4160

$ws.Range("A68").Value = "Result:"
$ws.Range("B68").Formula = '=IF(C66<=F66,"PASS","FAIL")'
$ws.Range("B68").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3) New "4. Visual inspection" section (rows 71, 73-77)
# ---------------------------------------------------------------------------

$ws.Range("A71").Value = "4. Visual inspection"
$ws.Range("A71").Font.Bold = $true
$ws.Range("A71").Font.Size = 14
$ws.Rows(71).RowHeight = 18

$ws.Range("A73").Value = "1) Jumper for LED installed"
$ws.Range("B73").Font.Bold = $true
$ws.Range("D73").Value = "Y"

$ws.Range("A74").Value = "2) Kapton foil present on back side"
$ws.Range("B74").Font.Bold = $true
$ws.Range("D74").Value = "Y"

$ws.Range("A75").Value = "3) S/N sticker attached"
$ws.Range("B75").Font.Bold = $true
$ws.Range("D75").Value = "Y"

$ws.Range("B76").Font.Bold = $true
$ws.Range("B77").Font.Bold = $true

# Row 75 (old "4. Summary" header, ht=18) and any other row dragged along by
# the insert/clear above should fall back to the default row height.
$ws.Rows(75).AutoFit()

# ---------------------------------------------------------------------------
# 4) Footer row (old row 72 -> new row 79): same Tester/Frank Meier/Test
#    date labels, only the date value itself moves forward by one day.
# ---------------------------------------------------------------------------

$ws.Range("F79").Value = 41830

# ---------------------------------------------------------------------------
# 5) Sheet-level bookkeeping: dimension grows, print area grows, selection
#    moves to the new last cell of interest.
# ---------------------------------------------------------------------------

$ws.PageSetup.PrintArea = '$A$1:$I$79'
$ws.Range("B79").Select()
